# daily auto push: 2026-02-20 14:04 UTC
#
# Inserts one new data row (2026/02/20, 金, 20, 201) immediately before the
# existing row 848 of Sheet1, pushing the rows that used to be 848..889 down
# to 849..890. The sheet's used-range grows from A1:D889 to A1:D890.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 848..889 down by one to make room for the new record.
$ws.Rows.Item(848).Insert()

# The date column holds plain text like "2026/12/29", not real dates, so
# force the new cell to Text formatting before writing the string -
# otherwise Excel would auto-convert "2026/02/20" into a date serial value.
$ws.Range("A848").NumberFormat = "@"
$ws.Range("A848").Value = "2026/02/20"
$ws.Range("B848").Value = "金"
$ws.Range("C848").Value = 20
$ws.Range("D848").Value = 201
